$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
  3  = 22.0
  4  = 33.0
  5  = 44.0
  6  = 55.0
  7  = 66.0
  8  = 77.0
  9  = 88.0
  10 = 99.0
  11 = 110.0
  12 = 121.0
  13 = 132.0
  14 = 143.0
  15 = 154.0
  16 = 165.0
  17 = 176.0
  18 = 187.0
  19 = 198.0
  20 = 209.0
  21 = 220.0
}

foreach ($row in $values.Keys) {
  $ws.Range("B$row").Value = $values[$row]
}
